$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5 ("第四天" row): fill in completion status, completion time, and notes
$ws.Range("C5").Value = "完成"
$ws.Range("D5").NumberFormat = $ws.Range("D4").NumberFormat
$ws.Range("D5").Value = 0.333333333333333
$ws.Range("F5").Value = "完成了合并以及ppt制作"

# Row 6 ("第五天" row): fill in planned content
$ws.Range("B6").Value = "项目答辩"

# Match the final selection shown in the diff
$ws.Range("F5").Select()
